$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that hold numeric-looking price strings must be forced to
# Text format first, otherwise Excel auto-converts the assigned string into a
# number (dropping e.g. trailing zeros: "0.450" -> 0.45).
$textCells = @('D4','D5','D6','D8','D10','D12','D13','D17','D20','D21','D22','D23','D24','D25','D26','D28','D31','D32','D33','D34','D35','D36','D37','D40','D41','D43','D45','D46','D48','D50','D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.092.19'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '3.794.38'
$ws.Range('D4').Value = '0.995'
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = '600.99'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = '164.95'
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.517'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').Value = '0.450'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  +2.88%  '
$ws.Range('D12').Value = '0.0000248'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('D13').Value = '35.67'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '4.432.20'
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = '3.802.91'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = '68.105.86'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '18.29'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '461.38'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = '9.70'
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').Value = '0.702'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '0.0000147'
$ws.Range('E23').Value = '  -4.24%  '
$ws.Range('D24').Value = '83.08'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').Value = '11.98'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('D26').Value = '2.11'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '9.99'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').Value = '3.944.79'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '7.34'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '2.63'
$ws.Range('E32').Value = '  -4.94%  '
$ws.Range('D33').Value = '29.30'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '9.02'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('D36').Value = '0.0998'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = '3.32'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('D40').Value = '0.988'
$ws.Range('E40').Value = '  -1.46%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D43').Value = '47.56'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').Value = '43.23'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = '152.30'
$ws.Range('E46').Value = '  +2.54%  '
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').Value = '1.87'
$ws.Range('E48').Value = '  +1.89%  '
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('D50').Value = '390.08'
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('D51').Value = '26.46'
$ws.Range('E51').Value = '  -1.48%  '
